$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new data row: B8 = 4
$ws.Range("B8").Value = 4

# Update the shared string text (remove trailing " !")
$ws.Range("A1").Value = "On inscrit le nombre d'heures régulièrement… Le contenu du travail n'est pas important"

# Move the active selection to D8
$ws.Range("D8").Select()
